# compte rendu_texte.docx - "version corrigée du compte rendu"
# Applies the tracked content changes; sentences that were modified are
# recoloured in red (FF0000), matching the author's commit message.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Cosmetic run merges (remove spell/grammar proofErr splits) that
#    carry no visible text change. Re-"typing" the same text over the
#    found span makes the host re-tokenize the paragraph into plain
#    runs, mirroring the diff's removal of <w:proofErr/> markers.
# ------------------------------------------------------------------

$t = "format dit bmp. Il s’agit alors d’enregistrer"
$null = $d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

$t = "Le format JEPG (Joint Photographic Expert Group) s’avère être un bon compromis entre une qualité d’image très correcte et une bonne compression. La compression JPEG s’appuie sur "
$null = $d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

$t = "la DCT (Discret Cosinus Transform)."
$null = $d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

$t = "Nous avons réalisé à l’aide du logicile de programmation Python une compression de type JPEG."
$null = $d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

$t = " suit le traitement suivant :"
$null = $d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

$t = "une RLE (Run-Lenth Encoding), qui consiste"
$null = $d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

# ------------------------------------------------------------------
# 2) Grammar fix: "un image finale" -> "une image finale"
# ------------------------------------------------------------------

$t = "pour donner un image finale de même taille"
$r = "pour donner une image finale de même taille"
$null = $d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $r, 2)

# ------------------------------------------------------------------
# 3) New paragraph inserted before "La Matrice ainsi obtenue..." :
#    "Nous avons retenu la matrice de quantification de la norme
#    JPEG, étudiée[_GoBack] pour un rendu optimal." in red.
#    The document's sole "_GoBack" bookmark moves here too.
# ------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute("La Matrice ainsi obtenue", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$null = $para.Range.InsertParagraphBefore()
$para.Range.Text = "Nous avons retenu la matrice de quantification de la norme JPEG, étudiée pour un rendu optimal."
$para.Range.Font.ItalicBi = $true
$para.Range.Font.Color = 255

$rng3 = $d.Content
$null = $rng3.Find.Execute("étudiée", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($rng3.End, $rng3.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 4) New sentence inserted before "Cette liste ..." (red):
#    "On remarque alors que les hautes « fréquences » se retrouvent
#    à la fin de la liste. "
# ------------------------------------------------------------------

$rng4 = $d.Content
$null = $rng4.Find.Execute("Cette liste ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertionPoint = $d.Range($rng4.Start, $rng4.Start)
$insertionPoint.InsertBefore("On remarque alors que les hautes « fréquences » se retrouvent à la fin de la liste. ")

$rng5 = $d.Content
$null = $rng5.Find.Execute("On remarque alors que les hautes « fréquences » se retrouvent à la fin de la liste. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng5.Font.ItalicBi = $true
$rng5.Font.Color = 255

$d.Save()
